# Applies the "sapt. 4" (week 4, column F) presence update for the
# students on rows 4, 6-14, 16-21: increment the attendance count in
# column F by 1 (empty cells become 1, cells holding 1 become 2).
# Column Q recalculates automatically since it's driven by a SUM formula.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    4  = 1
    6  = 2
    7  = 1
    8  = 1
    9  = 2
    10 = 2
    11 = 1
    12 = 2
    13 = 2
    14 = 2
    16 = 1
    17 = 2
    18 = 2
    19 = 2
    20 = 1
    21 = 2
}

foreach ($r in $newValues.Keys) {
    $ws.Cells.Item($r, 6).Value = $newValues[$r]   # column F is the 6th column
}

$wb.Save()
